$d = $word.ActiveDocument

# Map of exact paragraph text (as returned by Range.Text, without the
# trailing paragraph mark) to the new text that should replace it.
$replacements = @{
    "I helped with IAS interfaces and it was nice." = "Design: I helped with IAS interfaces and it was nice.";
    "Yes" = "Design: Yes";
    "The co-operation with Class LR was top professional and very well remembered." = "Design: The co-operation with Class LR was top professional and very well remembered.";
    "No problems" = "Design: No problems";
    "Our small purchases for mounting parts was made difficult by payment reputation of ours. Also the delivery lead times of many materials and parts was long." = "Design: Our small purchases for mounting parts was made difficult by payment reputation of ours. Also the delivery lead times of many materials and parts was long.";
    "Budget was good, yard was refunded." = "Design: Budget was good, yard was refunded.";
    "Was good." = "Design: Was good.";
    "First time in 20 years automation design made all drawings, not outsourced. " = "Design: First time in 20 years automation design made all drawings, not outsourced. ";
    "Pleasant and efficient co-operation in spite of the very heavy workload of the good electrical design colleagues." = "Design: Pleasant and efficient co-operation in spite of the very heavy workload of the good electrical design colleagues.";
    "Machinery design have skilled designers." = "Design: Machinery design have skilled designers.";
    "Was OK." = "Design: Was OK.";
    "Jira and ERM  doesn't work as I expected. ACAD without electrical symbols." = "Design: Jira and ERM  doesn't work as I expected. ACAD without electrical symbols.";
    "More teamwork before purchase of different sister systems. " = "Design: More teamwork before purchase of different sister systems. "
}

# Find the paragraph indexes whose visible text (paragraph mark stripped)
# matches one of the keys above.
$paragraphs = $d.Paragraphs
$matchedIndexes = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $text = $paragraphs.Item($i).Range.Text
    $trimmed = $text.TrimEnd([char]13, [char]7)
    if ($replacements.ContainsKey($trimmed)) {
        [void]$matchedIndexes.Add($i)
    }
}

# Apply the edits from the end of the document towards the start, and use
# a Range that spans only the paragraph's visible text (i.e. excludes the
# trailing paragraph/cell mark) so we don't insert an extra paragraph.
# Processing in reverse guarantees earlier (not-yet-processed) paragraphs'
# offsets are unaffected by edits made to later paragraphs.
for ($j = $matchedIndexes.Count - 1; $j -ge 0; $j--) {
    $i = $matchedIndexes[$j]
    $rng = $paragraphs.Item($i).Range
    $text = $rng.Text
    $trimmed = $text.TrimEnd([char]13, [char]7)
    $newRng = $d.Range($rng.Start, $rng.Start + $trimmed.Length)
    $newRng.Text = $replacements[$trimmed]
}
